$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = 1.14725446361153
$ws.Range("D2").Value = 0.2592856585496741

# Row 3
$ws.Range("C3").Value = 0.4060421473599262
$ws.Range("D3").Value = 0.6872559169804111

# Row 4
$ws.Range("C4").Value = 1.538821615351736
$ws.Range("D4").Value = 0.1331037356110538
$ws.Range("G4").Value = "No"

# Row 5
$ws.Range("C5").Value = 0.0171273209410233
$ws.Range("D5").Value = 0.98643515741041

# Row 6
$ws.Range("C6").Value = -1.120432360801508
$ws.Range("D6").Value = 0.2703817234616073

# Row 7
$ws.Range("C7").Value = 0.4666371027944691
$ws.Range("D7").Value = 0.6437358818395951

# Row 8
$ws.Range("C8").Value = -1.500201011762955
$ws.Range("D8").Value = 0.1427920965003817

# Row 9
$ws.Range("C9").Value = 1.377221673351219
$ws.Range("D9").Value = 0.1774437468296994

# Row 10
$ws.Range("C10").Value = -0.6921621030769461
$ws.Range("D10").Value = 0.4935358373528718

# Row 11
$ws.Range("C11").Value = -1.436415213468986
$ws.Range("D11").Value = 0.1600225433442253
